$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 185295
$ws.Cells.Item(2, 4).Value = 238956444
$ws.Cells.Item(6, 3).Value = 501
$ws.Cells.Item(6, 4).Value = 744414
$ws.Cells.Item(8, 3).Value = 76513
$ws.Cells.Item(8, 4).Value = 112899229
$ws.Cells.Item(10, 3).Value = 36118
$ws.Cells.Item(10, 4).Value = 52445295
$ws.Cells.Item(13, 3).Value = 2069
$ws.Cells.Item(13, 4).Value = 2948472
$ws.Cells.Item(16, 3).Value = 2380
$ws.Cells.Item(16, 4).Value = 3353634
$ws.Cells.Item(17, 3).Value = 47038
$ws.Cells.Item(17, 4).Value = 59832863
$ws.Cells.Item(19, 3).Value = 21
$ws.Cells.Item(19, 4).Value = 30563
$ws.Cells.Item(21, 3).Value = 169
$ws.Cells.Item(21, 4).Value = 247241
$ws.Cells.Item(22, 3).Value = 21669
$ws.Cells.Item(22, 4).Value = 31935155
$ws.Cells.Item(24, 3).Value = 7162
$ws.Cells.Item(24, 4).Value = 10391888
$ws.Cells.Item(26, 3).Value = 845
$ws.Cells.Item(26, 4).Value = 1188028
$ws.Cells.Item(28, 3).Value = 722
$ws.Cells.Item(28, 4).Value = 1023983
$ws.Cells.Item(29, 3).Value = 60009
$ws.Cells.Item(29, 4).Value = 77058224
$ws.Cells.Item(30, 3).Value = 32
$ws.Cells.Item(30, 4).Value = 37462
$ws.Cells.Item(32, 3).Value = 511
$ws.Cells.Item(32, 4).Value = 758971
$ws.Cells.Item(34, 3).Value = 29765
$ws.Cells.Item(34, 4).Value = 43921013
$ws.Cells.Item(36, 3).Value = 5661
$ws.Cells.Item(36, 4).Value = 8159694
$ws.Cells.Item(38, 3).Value = 761
$ws.Cells.Item(38, 4).Value = 1062945
$ws.Cells.Item(39, 3).Value = 791
$ws.Cells.Item(39, 4).Value = 1114671
$ws.Cells.Item(40, 3).Value = 41449
$ws.Cells.Item(40, 4).Value = 52652896
$ws.Cells.Item(41, 3).Value = 15
$ws.Cells.Item(41, 4).Value = 14595
$ws.Cells.Item(44, 3).Value = 205
$ws.Cells.Item(44, 4).Value = 303208
$ws.Cells.Item(45, 3).Value = 18661
$ws.Cells.Item(45, 4).Value = 27516106
$ws.Cells.Item(47, 3).Value = 6809
$ws.Cells.Item(47, 4).Value = 9883782
$ws.Cells.Item(48, 3).Value = 775
$ws.Cells.Item(48, 4).Value = 1086606
$ws.Cells.Item(50, 3).Value = 515
$ws.Cells.Item(50, 4).Value = 729983
$ws.Cells.Item(51, 3).Value = 11565
$ws.Cells.Item(51, 4).Value = 15154115
$ws.Cells.Item(55, 3).Value = 4280
$ws.Cells.Item(55, 4).Value = 6282849
$ws.Cells.Item(56, 3).Value = 2907
$ws.Cells.Item(56, 4).Value = 4236301
$ws.Cells.Item(57, 3).Value = 244
$ws.Cells.Item(57, 4).Value = 339625
$ws.Cells.Item(58, 3).Value = 106
$ws.Cells.Item(58, 4).Value = 146116
$ws.Cells.Item(59, 3).Value = 84258
$ws.Cells.Item(59, 4).Value = 106517955
$ws.Cells.Item(63, 3).Value = 234
$ws.Cells.Item(63, 4).Value = 341636
$ws.Cells.Item(65, 3).Value = 41983
$ws.Cells.Item(65, 4).Value = 61986095
$ws.Cells.Item(66, 3).Value = 47
$ws.Cells.Item(66, 4).Value = 69674
$ws.Cells.Item(68, 3).Value = 18578
$ws.Cells.Item(68, 4).Value = 27013388
$ws.Cells.Item(70, 3).Value = 1363
$ws.Cells.Item(70, 4).Value = 1962434
$ws.Cells.Item(71, 3).Value = 1038
$ws.Cells.Item(71, 4).Value = 1472682
$ws.Cells.Item(72, 3).Value = 13769
$ws.Cells.Item(72, 4).Value = 18660209
$ws.Cells.Item(76, 3).Value = 4303
$ws.Cells.Item(76, 4).Value = 6346768
$ws.Cells.Item(78, 3).Value = 3385
$ws.Cells.Item(78, 4).Value = 4887719
$ws.Cells.Item(80, 3).Value = 198
$ws.Cells.Item(80, 4).Value = 284311
$ws.Cells.Item(81, 3).Value = 159
$ws.Cells.Item(81, 4).Value = 230982
$ws.Cells.Item(82, 3).Value = 2915
$ws.Cells.Item(82, 4).Value = 4012786
$ws.Cells.Item(83, 3).Value = 862
$ws.Cells.Item(83, 4).Value = 1267254
$ws.Cells.Item(85, 3).Value = 1049
$ws.Cells.Item(85, 4).Value = 1524552
$ws.Cells.Item(88, 3).Value = 86277
$ws.Cells.Item(88, 4).Value = 107763342
$ws.Cells.Item(90, 3).Value = 40
$ws.Cells.Item(90, 4).Value = 59061
$ws.Cells.Item(92, 3).Value = 537
$ws.Cells.Item(92, 4).Value = 796277
$ws.Cells.Item(94, 3).Value = 35613
$ws.Cells.Item(94, 4).Value = 52491684
$ws.Cells.Item(96, 3).Value = 16760
$ws.Cells.Item(96, 4).Value = 24370536
$ws.Cells.Item(97, 3).Value = 669
$ws.Cells.Item(97, 4).Value = 910382
$ws.Cells.Item(98, 3).Value = 847
$ws.Cells.Item(98, 4).Value = 1191872
$ws.Cells.Item(100, 3).Value = 24768
$ws.Cells.Item(100, 4).Value = 33430461
$ws.Cells.Item(104, 3).Value = 316
$ws.Cells.Item(104, 4).Value = 474000
$ws.Cells.Item(105, 3).Value = 9047
$ws.Cells.Item(105, 4).Value = 13340314
$ws.Cells.Item(106, 3).Value = 2304
$ws.Cells.Item(106, 4).Value = 3338722
$ws.Cells.Item(108, 3).Value = 208
$ws.Cells.Item(108, 4).Value = 300215
$ws.Cells.Item(109, 3).Value = 157
$ws.Cells.Item(109, 4).Value = 217288
$ws.Cells.Item(110, 3).Value = 9062
$ws.Cells.Item(110, 4).Value = 12078830
$ws.Cells.Item(112, 3).Value = 4156
$ws.Cells.Item(112, 4).Value = 6082127
$ws.Cells.Item(113, 3).Value = 2599
$ws.Cells.Item(113, 4).Value = 3757705
$ws.Cells.Item(114, 3).Value = 129
$ws.Cells.Item(114, 4).Value = 180513
$ws.Cells.Item(116, 3).Value = 3045
$ws.Cells.Item(116, 4).Value = 4313809
$ws.Cells.Item(117, 3).Value = 775
$ws.Cells.Item(117, 4).Value = 1148943
$ws.Cells.Item(118, 3).Value = 106
$ws.Cells.Item(118, 4).Value = 156600
$ws.Cells.Item(120, 3).Value = 18
$ws.Cells.Item(120, 4).Value = 26949
$ws.Cells.Item(121, 3).Value = 54210
$ws.Cells.Item(121, 4).Value = 68863412
$ws.Cells.Item(125, 3).Value = 364
$ws.Cells.Item(125, 4).Value = 536762
$ws.Cells.Item(127, 3).Value = 22993
$ws.Cells.Item(127, 4).Value = 33905865
$ws.Cells.Item(129, 3).Value = 8095
$ws.Cells.Item(129, 4).Value = 11753126
$ws.Cells.Item(131, 3).Value = 693
$ws.Cells.Item(131, 4).Value = 976931
$ws.Cells.Item(132, 3).Value = 673
$ws.Cells.Item(132, 4).Value = 959426
$ws.Cells.Item(133, 3).Value = 143767
$ws.Cells.Item(133, 4).Value = 181133519
$ws.Cells.Item(134, 3).Value = 64
$ws.Cells.Item(134, 4).Value = 68835
$ws.Cells.Item(138, 3).Value = 518
$ws.Cells.Item(138, 4).Value = 762059
$ws.Cells.Item(139, 3).Value = 6
$ws.Cells.Item(139, 4).Value = 8383
$ws.Cells.Item(140, 3).Value = 58174
$ws.Cells.Item(140, 4).Value = 85661375
$ws.Cells.Item(143, 3).Value = 20639
$ws.Cells.Item(143, 4).Value = 29846607
$ws.Cells.Item(146, 3).Value = 2817
$ws.Cells.Item(146, 4).Value = 4029106
$ws.Cells.Item(148, 3).Value = 1918
$ws.Cells.Item(148, 4).Value = 2667013
$ws.Cells.Item(149, 3).Value = 4
$ws.Cells.Item(149, 4).Value = 5215
$ws.Cells.Item(150, 3).Value = 154413
$ws.Cells.Item(150, 4).Value = 192931258
$ws.Cells.Item(151, 3).Value = 79
$ws.Cells.Item(151, 4).Value = 78978
$ws.Cells.Item(155, 3).Value = 356
$ws.Cells.Item(155, 4).Value = 524558
$ws.Cells.Item(157, 3).Value = 61970
$ws.Cells.Item(157, 4).Value = 91147484
$ws.Cells.Item(160, 3).Value = 31443
$ws.Cells.Item(160, 4).Value = 45658010
$ws.Cells.Item(163, 3).Value = 2532
$ws.Cells.Item(163, 4).Value = 3572234
$ws.Cells.Item(166, 3).Value = 2213
$ws.Cells.Item(166, 4).Value = 3087613
$ws.Cells.Item(169, 3).Value = 63961
$ws.Cells.Item(169, 4).Value = 81452743
$ws.Cells.Item(173, 3).Value = 313
$ws.Cells.Item(173, 4).Value = 456629
$ws.Cells.Item(175, 3).Value = 33397
$ws.Cells.Item(175, 4).Value = 49223407
$ws.Cells.Item(177, 3).Value = 7552
$ws.Cells.Item(177, 4).Value = 10897644
$ws.Cells.Item(179, 3).Value = 1145
$ws.Cells.Item(179, 4).Value = 1649505
$ws.Cells.Item(181, 3).Value = 1020
$ws.Cells.Item(181, 4).Value = 1425577
$ws.Cells.Item(182, 3).Value = 152302
$ws.Cells.Item(182, 4).Value = 193821304
$ws.Cells.Item(188, 3).Value = 491
$ws.Cells.Item(188, 4).Value = 729450
$ws.Cells.Item(190, 3).Value = 62297
$ws.Cells.Item(190, 4).Value = 91799253
$ws.Cells.Item(192, 3).Value = 39510
$ws.Cells.Item(192, 4).Value = 57466904
$ws.Cells.Item(194, 3).Value = 1345
$ws.Cells.Item(194, 4).Value = 1916481
$ws.Cells.Item(196, 3).Value = 1825
$ws.Cells.Item(196, 4).Value = 2560669
$ws.Cells.Item(197, 3).Value = 238382
$ws.Cells.Item(197, 4).Value = 314501238
$ws.Cells.Item(202, 3).Value = 754
$ws.Cells.Item(202, 4).Value = 1123218
$ws.Cells.Item(204, 3).Value = 125019
$ws.Cells.Item(204, 4).Value = 184669180
$ws.Cells.Item(205, 3).Value = 217
$ws.Cells.Item(205, 4).Value = 323380
$ws.Cells.Item(207, 3).Value = 99414
$ws.Cells.Item(207, 4).Value = 144909790
$ws.Cells.Item(209, 3).Value = 1250
$ws.Cells.Item(209, 4).Value = 1763849
$ws.Cells.Item(211, 3).Value = 2336
$ws.Cells.Item(211, 4).Value = 3310012
$ws.Cells.Item(212, 3).Value = 11
$ws.Cells.Item(212, 4).Value = 16500
